$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) for the two new date columns ("25-jul", "26-jul")
$ws.Range("AH1").Value = "25-jul"
$ws.Range("AI1").Value = "26-jul"

# Copy the formatting (number format + centered alignment) already used by the
# existing date-value columns onto the two new columns before writing values,
# so the new cells pick up the same style the rest of the table uses.
$ws.Range("AG2:AG11").Copy() | Out-Null
$ws.Range("AH2:AH11").PasteSpecial(-4122) | Out-Null
$ws.Range("AG2:AG11").Copy() | Out-Null
$ws.Range("AI2:AI11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data values for rows 2-11 (AH = "25-jul" column, AI = "26-jul" column)
$values = @{
    2  = @(11, 11)
    3  = @(16, 17)
    4  = @(9, 10)
    5  = @(16, 14)
    6  = @(18, 17)
    7  = @(12, 13)
    8  = @(17, 15)
    9  = @(16, 15)
    10 = @(23, 22)
    11 = @(19, 18)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 34).Value = $pair[0]   # column AH = 34
    $ws.Cells.Item($row, 35).Value = $pair[1]   # column AI = 35
}

# Leave the selection where the author's cursor ended up after data entry.
$ws.Range("AJ8").Select() | Out-Null
